$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.1581976666666667
$ws.Range("H2").Value = 0.474593
$ws.Range("I2").Value = 0.1400666049254827
$ws.Range("J2").Value = 0.1400666049254826
$ws.Range("M2").Value = 168.1098273333333
$ws.Range("N2").Value = 504.329482
$ws.Range("O2").Value = 0.2984182258032519
$ws.Range("P2").Value = 0.298418225803252
$ws.Range("Q2").Value = 26.59458242786956
$ws.Range("R2").Value = 239.351241850826
$ws.Range("S2").Value = 0.04179842773614757
$ws.Range("T2").Value = 0.04179842773614757
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.1581976666666667
$ws.Range("H3").Value = 0.474593
$ws.Range("I3").Value = 0.1400666049254827
$ws.Range("J3").Value = 0.1400666049254826
$ws.Range("O3").Value = 0.2893586437755394
$ws.Range("P3").Value = 0.2893586437755394
$ws.Range("Q3").Value = 25.78720613458345
$ws.Range("R3").Value = 232.084855211251
$ws.Range("S3").Value = 0.04052948283948196
$ws.Range("T3").Value = 0.04052948283948195
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.1581976666666667
$ws.Range("H4").Value = 0.474593
$ws.Range("I4").Value = 0.1400666049254827
$ws.Range("J4").Value = 0.1400666049254826
$ws.Range("M4").Value = 165.99353
$ws.Range("N4").Value = 497.98059
$ws.Range("O4").Value = 0.294661504941043
$ws.Range("P4").Value = 0.294661504941043
$ws.Range("Q4").Value = 26.25978912776333
$ws.Range("R4").Value = 236.33810214987
$ws.Range("S4").Value = 0.04127223659932523
$ws.Range("T4").Value = 0.04127223659932522
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 0.1581976666666667
$ws.Range("H5").Value = 0.474593
$ws.Range("I5").Value = 0.1400666049254827
$ws.Range("J5").Value = 0.1400666049254826
$ws.Range("M5").Value = 66.22673433333334
$ws.Range("N5").Value = 198.680203
$ws.Range("O5").Value = 0.1175616254801657
$ws.Range("P5").Value = 0.1175616254801657
$ws.Range("Q5").Value = 10.47691484248656
$ws.Range("R5").Value = 94.29223358237901
$ws.Range("S5").Value = 0.01646645775052792
$ws.Range("T5").Value = 0.01646645775052792
$ws.Range("I6").Value = 0.7029419733214338
$ws.Range("J6").Value = 0.7029419733214337
$ws.Range("M6").Value = 168.1098273333333
$ws.Range("N6").Value = 504.329482
$ws.Range("O6").Value = 0.2984182258032519
$ws.Range("P6").Value = 0.298418225803252
$ws.Range("Q6").Value = 133.46827576389
$ws.Range("R6").Value = 1201.21448187501
$ws.Range("S6").Value = 0.2097706965212192
$ws.Range("T6").Value = 0.2097706965212192
$ws.Range("I7").Value = 0.7029419733214338
$ws.Range("J7").Value = 0.7029419733214337
$ws.Range("O7").Value = 0.2893586437755394
$ws.Range("P7").Value = 0.2893586437755394
$ws.Range("S7").Value = 0.2034023360531915
$ws.Range("T7").Value = 0.2034023360531915
$ws.Range("I8").Value = 0.7029419733214338
$ws.Range("J8").Value = 0.7029419733214337
$ws.Range("M8").Value = 165.99353
$ws.Range("N8").Value = 497.98059
$ws.Range("O8").Value = 0.294661504941043
$ws.Range("P8").Value = 0.294661504941043
$ws.Range("Q8").Value = 131.78807324055
$ws.Range("R8").Value = 1186.09265916495
$ws.Range("S8").Value = 0.2071299397451202
$ws.Range("T8").Value = 0.2071299397451201
$ws.Range("I9").Value = 0.7029419733214338
$ws.Range("J9").Value = 0.7029419733214337
$ws.Range("M9").Value = 66.22673433333334
$ws.Range("N9").Value = 198.680203
$ws.Range("O9").Value = 0.1175616254801657
$ws.Range("P9").Value = 0.1175616254801657
$ws.Range("Q9").Value = 52.57972232293501
$ws.Range("R9").Value = 473.217500906415
$ws.Range("S9").Value = 0.08263900100190301
$ws.Range("T9").Value = 0.082639001001903
$ws.Range("G10").Value = 0.1199896666666667
$ws.Range("H10").Value = 0.359969
$ws.Range("I10").Value = 0.106237630366274
$ws.Range("J10").Value = 0.106237630366274
$ws.Range("M10").Value = 168.1098273333333
$ws.Range("N10").Value = 504.329482
$ws.Range("O10").Value = 0.2984182258032519
$ws.Range("P10").Value = 0.298418225803252
$ws.Range("Q10").Value = 20.17144214511755
$ws.Range("R10").Value = 181.542979306058
$ws.Range("S10").Value = 0.03170324516744517
$ws.Range("T10").Value = 0.03170324516744516
$ws.Range("G11").Value = 0.1199896666666667
$ws.Range("H11").Value = 0.359969
$ws.Range("I11").Value = 0.106237630366274
$ws.Range("J11").Value = 0.106237630366274
$ws.Range("O11").Value = 0.2893586437755394
$ws.Range("P11").Value = 0.2893586437755394
$ws.Range("Q11").Value = 19.55906388223145
$ws.Range("R11").Value = 176.031574940083
$ws.Range("S11").Value = 0.03074077664071211
$ws.Range("T11").Value = 0.0307407766407121
$ws.Range("G12").Value = 0.1199896666666667
$ws.Range("H12").Value = 0.359969
$ws.Range("I12").Value = 0.106237630366274
$ws.Range("J12").Value = 0.106237630366274
$ws.Range("M12").Value = 165.99353
$ws.Range("N12").Value = 497.98059
$ws.Range("O12").Value = 0.294661504941043
$ws.Range("P12").Value = 0.294661504941043
$ws.Range("Q12").Value = 19.91750833352333
$ws.Range("R12").Value = 179.25757500171
$ws.Range("S12").Value = 0.03130414004509654
$ws.Range("T12").Value = 0.03130414004509653
$ws.Range("G13").Value = 0.1199896666666667
$ws.Range("H13").Value = 0.359969
$ws.Range("I13").Value = 0.106237630366274
$ws.Range("J13").Value = 0.106237630366274
$ws.Range("M13").Value = 66.22673433333334
$ws.Range("N13").Value = 198.680203
$ws.Range("O13").Value = 0.1175616254801657
$ws.Range("P13").Value = 0.1175616254801657
$ws.Range("Q13").Value = 7.946523777078556
$ws.Range("R13").Value = 71.518713993707
$ws.Range("S13").Value = 0.01248946851302018
$ws.Range("T13").Value = 0.01248946851302017
$ws.Range("E14").Value = 2
$ws.Range("F14").Value = 0.6666666666666666
$ws.Range("G14").Value = 0.05732366666666666
$ws.Range("H14").Value = 0.171971
$ws.Range("I14").Value = 0.05075379138680971
$ws.Range("J14").Value = 0.05075379138680969
$ws.Range("M14").Value = 168.1098273333333
$ws.Range("N14").Value = 504.329482
$ws.Range("O14").Value = 0.2984182258032519
$ws.Range("P14").Value = 0.298418225803252
$ws.Range("Q14").Value = 9.636671705446888
$ws.Range("R14").Value = 86.73004534902199
$ws.Range("S14").Value = 0.01514585637844012
$ws.Range("T14").Value = 0.01514585637844012
$ws.Range("E15").Value = 2
$ws.Range("F15").Value = 0.6666666666666666
$ws.Range("G15").Value = 0.05732366666666666
$ws.Range("H15").Value = 0.171971
$ws.Range("I15").Value = 0.05075379138680971
$ws.Range("J15").Value = 0.05075379138680969
$ws.Range("O15").Value = 0.2893586437755394
$ws.Range("P15").Value = 0.2893586437755394
$ws.Range("Q15").Value = 9.34411511794411
$ws.Range("R15").Value = 84.097036061497
$ws.Range("S15").Value = 0.01468604824215391
$ws.Range("T15").Value = 0.01468604824215391
$ws.Range("E16").Value = 2
$ws.Range("F16").Value = 0.6666666666666666
$ws.Range("G16").Value = 0.05732366666666666
$ws.Range("H16").Value = 0.171971
$ws.Range("I16").Value = 0.05075379138680971
$ws.Range("J16").Value = 0.05075379138680969
$ws.Range("M16").Value = 165.99353
$ws.Range("N16").Value = 497.98059
$ws.Range("O16").Value = 0.294661504941043
$ws.Range("P16").Value = 0.294661504941043
$ws.Range("Q16").Value = 9.515357782543331
$ws.Range("R16").Value = 85.63822004289
$ws.Range("S16").Value = 0.01495518855150109
$ws.Range("T16").Value = 0.01495518855150109
$ws.Range("E17").Value = 2
$ws.Range("F17").Value = 0.6666666666666666
$ws.Range("G17").Value = 0.05732366666666666
$ws.Range("H17").Value = 0.171971
$ws.Range("I17").Value = 0.05075379138680971
$ws.Range("J17").Value = 0.05075379138680969
$ws.Range("M17").Value = 66.22673433333334
$ws.Range("N17").Value = 198.680203
$ws.Range("O17").Value = 0.1175616254801657
$ws.Range("P17").Value = 0.1175616254801657
$ws.Range("Q17").Value = 3.796359243345889
$ws.Range("R17").Value = 34.167233190113
$ws.Range("S17").Value = 0.005966698214714581
$ws.Range("T17").Value = 0.005966698214714579
Write-Host "Updated NATMI LR-pair values with new TPM data"
